# AGREE Evaluation aggregate - update aggregate data
# - Round columns B:E (rows 2-10) to 1 decimal place (values replaced)
# - Apply a "0.0" number format to columns F:G (rows 2-10), values unchanged
# - Append two new guideline rows (11: TRIPOD-AI, 12: APPRAISE-AI) with their scores
# - Move the active selection to L10 (matches the author's last on-screen selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 2-10: columns B-E get their stored value rounded to 1 dp ---
$existingBE = @{
    2  = @(88.9, 70.8, 70.3, 70.8)
    3  = @(87.5, 83.3, 76, 72.2)
    4  = @(73.6, 40.3, 53.1, 56.9)
    5  = @(91.6, 84.7, 75.5, 81.9)
    6  = @(68.1, 56.9, 35.4, 43.1)
    7  = @(90.3, 72.2, 66.2, 70.8)
    8  = @(93.3, 90.3, 79.2, 83.3)
    9  = @(91.7, 44.4, 31.3, 75)
    10 = @(93.1, 69.4, 68.2, 76.4)
}

foreach ($row in $existingBE.Keys) {
    $vals = $existingBE[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
}

# --- New rows 11-12: Guideline name (A) + scores for columns B-G ---
$ws.Cells.Item(11, 1).Value = "TRIPOD-AI"
$ws.Cells.Item(11, 2).Value = 84.7
$ws.Cells.Item(11, 3).Value = 93.1
$ws.Cells.Item(11, 4).Value = 79.5
$ws.Cells.Item(11, 5).Value = 75
$ws.Cells.Item(11, 6).Value = 73.6
$ws.Cells.Item(11, 7).Value = 91.7

$ws.Cells.Item(12, 1).Value = "APPRAISE-AI"
$ws.Cells.Item(12, 2).Value = 87.5
$ws.Cells.Item(12, 3).Value = 73.6
$ws.Cells.Item(12, 4).Value = 76
$ws.Cells.Item(12, 5).Value = 69.8
$ws.Cells.Item(12, 6).Value = 66.7
$ws.Cells.Item(12, 7).Value = 89.6

# --- Columns F:G, rows 2-12: apply a one-decimal display format (0.0) ---
# (existing F/G values for rows 2-10 are left as-is; only the display format changes)
$ws.Range("F2:G12").NumberFormat = "0.0"

# --- Match the workbook's last on-screen selection ---
$ws.Range("L10").Select() | Out-Null
